$wb = $excel.ActiveWorkbook

# --- table_specific_translations sheet: rename "individual" rows to "member" ---
$tst = $wb.Worksheets.Item("table_specific_translations")

$tst.Range("A5").Value = "member_id"
$tst.Range("B5").Value = "Member ID"

$tst.Range("A6").Value = "custom_member_form_id"
$tst.Range("B6").Value = "Custom Member Form ID"

$tst.Range("A7").Value = "custom_member_row_id"
$tst.Range("B7").Value = "Custom Member Row ID"

# widen column A on this sheet
$tst.Columns.Item(1).ColumnWidth = 42

# --- switch the active/selected sheet from "model" to "table_specific_translations" ---
$tst.Activate()
$tst.Range("A14").Select()
